$d = $word.ActiveDocument
$d.Content.Find.Execute("DoorObjects.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ExitObjects.", 2)
